$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5213562
$ws.Range("J17").Value = 5958243.5
$ws.Range("L17").Value = 17874730.5
$ws.Range("N17").Value = -17875066.5

$ws.Range("H112").Value = 1200.7693
$ws.Range("J112").Value = 1346.3636
$ws.Range("L112").Value = 4039.0908
$ws.Range("N112").Value = -6255.0908

$ws.Range("H113").Value = 2412.5
$ws.Range("I113").Value = 1753.125
$ws.Range("J113").Value = 2940
$ws.Range("K113").Value = 1753.125
$ws.Range("L113").Value = 2940
$ws.Range("M113").Value = 1500.875
$ws.Range("N113").Value = -9448

$ws.Range("H127").Value = 987.0769
$ws.Range("I127").Value = 781.7778
$ws.Range("J127").Value = 1449
$ws.Range("K127").Value = 2345.3334
$ws.Range("L127").Value = 4347
$ws.Range("M127").Value = 2614.6666
$ws.Range("N127").Value = -14267

$ws.Range("H132").Value = 3773.0747
$ws.Range("I132").Value = 3690.8728
$ws.Range("J132").Value = 4149.8335
$ws.Range("K132").Value = 11072.6184
$ws.Range("L132").Value = 12449.5005
$ws.Range("M132").Value = -8542.618399999999
$ws.Range("N132").Value = -17509.5005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1085.6786
$ws.Range("I2").Value = 945.7059
$ws.Range("J2").Value = 1302
$ws.Range("K2").Value = 945.7059
$ws.Range("L2").Value = 1302
$ws.Range("M2").Value = -832.7059
$ws.Range("N2").Value = -1528

$ws.Range("H32").Value = 10397.11
$ws.Range("I32").Value = 5598.757
$ws.Range("J32").Value = 24053.96
$ws.Range("K32").Value = 5598.757
$ws.Range("L32").Value = 24053.96
$ws.Range("M32").Value = -5311.757
$ws.Range("N32").Value = -24627.96

$ws.Range("H61").Value = 2616.8948
$ws.Range("I61").Value = 1975
$ws.Range("J61").Value = 3499.5
$ws.Range("K61").Value = 1975
$ws.Range("L61").Value = 3499.5
$ws.Range("M61").Value = -1763
$ws.Range("N61").Value = -3923.5

$ws.Range("H97").Value = 2093.2856
$ws.Range("I97").Value = 2147.95
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 2147.95
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -1651.95
$ws.Range("N97").Value = -1992

$ws.Range("H102").Value = 1839.4
$ws.Range("I102").Value = 1839.4
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1839.4
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -217.4000000000001
$ws.Range("N102").ClearContents()

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H105").Value = 38000
$ws.Range("J105").Value = 38000
$ws.Range("L105").Value = 38000
$ws.Range("N105").Value = -44988

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H116").Value = 1085.6786
$ws.Range("I116").Value = 945.7059
$ws.Range("J116").Value = 1302
$ws.Range("K116").Value = 945.7059
$ws.Range("L116").Value = 1302
$ws.Range("M116").Value = 1348.2941
$ws.Range("N116").Value = -5890

$ws.Range("H136").Value = 2616.8948
$ws.Range("I136").Value = 1975
$ws.Range("J136").Value = 3499.5
$ws.Range("K136").Value = 5925
$ws.Range("L136").Value = 10498.5
$ws.Range("M136").Value = -3375
$ws.Range("N136").Value = -15598.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1085.6786
$ws.Range("I3").Value = 945.7059
$ws.Range("J3").Value = 1302
$ws.Range("K3").Value = 945.7059
$ws.Range("L3").Value = 1302
$ws.Range("M3").Value = -831.7059
$ws.Range("N3").Value = -1530

$ws.Range("H20").Value = 1255.8823
$ws.Range("I20").Value = 905
$ws.Range("J20").Value = 1757.1428
$ws.Range("K20").Value = 905
$ws.Range("L20").Value = 1757.1428
$ws.Range("M20").Value = -658
$ws.Range("N20").Value = -2251.1428

$ws.Range("H94").Value = 11633.263
$ws.Range("I94").Value = 1000
$ws.Range("J94").Value = 29861.715
$ws.Range("K94").Value = 1000
$ws.Range("L94").Value = 29861.715
$ws.Range("M94").Value = -549
$ws.Range("N94").Value = -30763.715

$ws.Range("H103").Value = 14552.333
$ws.Range("J103").Value = 14552.333
$ws.Range("L103").Value = 14552.333
$ws.Range("N103").Value = -16896.333

$ws.Range("H110").Value = 43500
$ws.Range("J110").Value = 43500
$ws.Range("L110").Value = 43500
$ws.Range("N110").Value = -51680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1600.96
$ws.Range("I31").Value = 962.2308
$ws.Range("J31").Value = 3865.5454
$ws.Range("K31").Value = 962.2308
$ws.Range("L31").Value = 3865.5454
$ws.Range("M31").Value = -667.2308
$ws.Range("N31").Value = -4455.5454

$ws.Range("H34").Value = 1600.96
$ws.Range("I34").Value = 962.2308
$ws.Range("J34").Value = 3865.5454
$ws.Range("K34").Value = 962.2308
$ws.Range("L34").Value = 3865.5454
$ws.Range("M34").Value = -760.2308
$ws.Range("N34").Value = -4269.5454

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H132").Value = 2366.7896
$ws.Range("I132").Value = 2677
$ws.Range("J132").Value = 1940.25
$ws.Range("K132").Value = 8031
$ws.Range("L132").Value = 5820.75
$ws.Range("M132").Value = -5501
$ws.Range("N132").Value = -10880.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 804.5333000000001
$ws.Range("J131").Value = 1090.8334
$ws.Range("L131").Value = 3272.5002
$ws.Range("N131").Value = -13352.5002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3126.2
$ws.Range("I61").Value = 2955.4119
$ws.Range("J61").Value = 4094
$ws.Range("K61").Value = 2955.4119
$ws.Range("L61").Value = 4094
$ws.Range("M61").Value = -2753.4119
$ws.Range("N61").Value = -4498

$ws.Range("H113").Value = 3126.2
$ws.Range("I113").Value = 2955.4119
$ws.Range("J113").Value = 4094
$ws.Range("K113").Value = 2955.4119
$ws.Range("L113").Value = 4094
$ws.Range("M113").Value = -785.4119000000001
$ws.Range("N113").Value = -8434
